$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.196.91'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").Value = '3.798.00'
$ws.Range("E3").Value = '  +1.46%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.42'
$ws.Range("E5").Value = '  -0.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '164.37'
$ws.Range("E6").Value = '  -2.63%  '
$ws.Range("D7").Value = '3.796.22'
$ws.Range("E7").Value = '  +1.38%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.538'
$ws.Range("E9").Value = '  +0.54%  '
$ws.Range("E10").Value = '  +2.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.32'
$ws.Range("E11").Value = '  -0.34%  '
$ws.Range("E12").Value = '  -0.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.38'
$ws.Range("E13").Value = '  -2.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000247'
$ws.Range("E14").Value = '  -0.62%  '
$ws.Range("D15").Value = '4.433.39'
$ws.Range("E15").Value = '  +1.47%  '
$ws.Range("D16").Value = '3.793.39'
$ws.Range("E16").Value = '  +1.38%  '
$ws.Range("D17").Value = '69.316.33'
$ws.Range("E17").Value = '  +0.62%  '
$ws.Range("E18").Value = '  +2.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.38'
$ws.Range("E19").Value = '  +1.62%  '
$ws.Range("E20").Value = '  -0.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.26'
$ws.Range("E21").Value = '  +3.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '491.84'
$ws.Range("E22").Value = '  -0.40%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.723'
$ws.Range("E23").Value = '  -0.53%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000153'
$ws.Range("E24").Value = '  -1.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.69'
$ws.Range("E25").Value = '  -0.68%  '
$ws.Range("E26").Value = '  -2.84%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.24'
$ws.Range("E27").Value = '  -1.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.11'
$ws.Range("E28").Value = '  -2.83%  '
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.10'
$ws.Range("E31").Value = '  +1.57%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.41'
$ws.Range("E32").Value = '  -4.76%  '
$ws.Range("D33").Value = '3.942.16'
$ws.Range("E33").Value = '  +1.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '31.95'
$ws.Range("E34").Value = '  +0.56%  '
$ws.Range("D35").Value = '3.746.39'
$ws.Range("E35").Value = '  +1.88%  '
$ws.Range("E36").Value = '  -1.43%  '
$ws.Range("E37").Value = '  +6.30%  '
$ws.Range("E38").Value = '  +0.23%  '
$ws.Range("E39").Value = '  +1.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  -0.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.324'
$ws.Range("E41").Value = '  -0.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.05'
$ws.Range("E42").Value = '  +1.73%  '
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.99'
$ws.Range("E43").Value = '  +1.07%  '
$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '425.70'
$ws.Range("E44").Value = '  -2.45%  '
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '48.43'
$ws.Range("E45").Value = '  -0.88%  '
$ws.Range("E47").Value = '  -0.69%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '142.19'
$ws.Range("E48").Value = '  +0.50%  '
$ws.Range("D49").Value = '2.829.69'
$ws.Range("E49").Value = '  +2.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '39.80'
$ws.Range("E50").Value = '  -1.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.29'
$ws.Range("E51").Value = '  +5.92%  '
